# "Sistemato il thread di invio dei messaggi"
# The D column held a mix of malformed / inconsistent phone numbers
# (+393466296727, +393802857952, +39 a 789, +3270386536, 33350364ice, ...
# plus a couple of blank rows). They all get consolidated to a single,
# correctly formatted number, and the two blank rows in the middle of the
# list (D3 and D6) are filled in with it as well so the whole D1:D7 block
# is uniform.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$phone = "+39 346 629 6727"

$ws.Range("D1").Value = $phone
$ws.Range("D2").Value = $phone
$ws.Range("D3").Value = $phone
$ws.Range("D4").Value = $phone
$ws.Range("D5").Value = $phone
$ws.Range("D6").Value = $phone
$ws.Range("D7").Value = $phone

# Move the active selection off D5 (where it was left before) to C16,
# matching where the author clicked next.
$null = $ws.Range("C16").Select()
